$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A23").Value = "8/29"
$ws.Range("B23").Value = "10/24"
$ws.Range("C23").Value = "第73期 秘寶 開放區域 水晶迷城 祕寶效果: 收穫騎乘獸門票有3%(10.8)機率翻倍"

$ws.Range("C24").Select()
